# Update cryptos list values per the scraped diff (commit message:
# "Updated cryptos list on Thu Aug 24 13:44:03 UTC 2023 with GitHub Actions").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '''26.445.15'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +2.15%  '

# Row 3
$ws.Range("D3").Value = '''1.669.99'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.80%  '

# Row 5
$ws.Range("D5").Value = '''220.24'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.19%  '

# Row 6
$ws.Range("D6").Value = '''0.5255'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.89%  '

# Row 7
$ws.Range("E7").Value = '  +0.22%  '

# Row 8
$ws.Range("D8").Value = '''0.2668'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +2.52%  '

# Row 9
$ws.Range("D9").Value = '''0.06364'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.70%  '

# Row 10
$ws.Range("D10").Value = '''21.66'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +5.24%  '

# Row 11
$ws.Range("D11").Value = '''0.07800'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.76%  '

# Row 12
$ws.Range("B12").Value = 'WrappedEther'
$ws.Range("C12").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D12").Value = '''1.676.23'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +2.19%  '

# Row 13
$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").Value = '''4.465'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.26%  '

# Row 14
$ws.Range("D14").Value = '''0.5527'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.91%  '

# Row 15
$ws.Range("D15").Value = '''0.0₅8268'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.09%  '

# Row 16
$ws.Range("D16").Value = '''65.45'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.60%  '

# Row 17
$ws.Range("D17").Value = '''26.458.86'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +2.13%  '

# Row 19
$ws.Range("D19").Value = '''4.739'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.14%  '

# Row 20
$ws.Range("D20").Value = '''193.49'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.76%  '

# Row 21
$ws.Range("D21").Value = '''10.33'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.85%  '

# Row 22
$ws.Range("D22").Value = '''6.264'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.40%  '

# Row 24
$ws.Range("B24").Value = 'Stellar'
$ws.Range("C24").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D24").Value = '''0.1259'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.85%  '

# Row 25
$ws.Range("B25").Value = 'Monero'
$ws.Range("C25").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D25").Value = '''138.79'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.91%  '

# Row 26
$ws.Range("D26").Value = '''7.393'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.51%  '

# Row 27
$ws.Range("E27").Value = '  +2.37%  '

# Row 28
$ws.Range("D28").Value = '''1.419'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.83%  '

# Row 29
$ws.Range("D29").Value = '''0.06135'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +3.51%  '

# Row 30
$ws.Range("D30").Value = '''1.290'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.61%  '

# Row 31
$ws.Range("D31").Value = '''3.611'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +6.30%  '

# Row 32
$ws.Range("D32").Value = '''3.394'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.18%  '

# Row 33
$ws.Range("D33").Value = '''1.683'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.97%  '

# Row 34
$ws.Range("E34").Value = '  +1.72%  '

# Row 35
$ws.Range("D35").Value = '''0.6067'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +8.39%  '

# Row 36
$ws.Range("D36").Value = '''2.422'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.09%  '

# Row 37
$ws.Range("E37").Value = '  +1.01%  '

# Row 38
$ws.Range("E38").Value = '  +0.78%  '

# Row 39
$ws.Range("D39").Value = '''6.031'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +3.65%  '

# Row 40
$ws.Range("D40").Value = '''1.090.99'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +7.01%  '

# Row 41
$ws.Range("D41").Value = '''0.8598'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.02%  '

# Row 42
$ws.Range("E42").Value = '  +0.11%  '

# Row 43
$ws.Range("D43").Value = '''100.68'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.20%  '

# Row 44
$ws.Range("D44").Value = '''1.812.08'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.30%  '

# Row 45
$ws.Range("E45").Value = '  +4.68%  '

# Row 46
$ws.Range("D46").Value = '''0.0₈109'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.65%  '

# Row 47
$ws.Range("D47").Value = '''8.170'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.47%  '

# Row 48
$ws.Range("E48").Value = '  +0.06%  '

# Row 49
$ws.Range("D49").Value = '''0.05205'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.31%  '

# Row 50
$ws.Range("D50").Value = '''1.485'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +8.64%  '

# Row 51
$ws.Range("D51").Value = '''0.4232'
$ws.Range("D51").Style = "Normal"
